# Refresh crypto price (D) and 1h volume-change (E) columns, plus the
# three newly-added/renamed coins in rows 49-51 (B/C/D/E), per the
# GitHub Actions data-refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''51.979.69'
$ws.Range("E2").Value = '''  +0.32%  '
$ws.Range("D3").Value = '''3.000.94'
$ws.Range("E3").Value = '''  +2.59%  '
$ws.Range("E4").Value = '''  -0.05%  '
$ws.Range("D5").Value = '''354.33'
$ws.Range("E5").Value = '''  -0.19%  '
$ws.Range("D6").Value = '''106.48'
$ws.Range("E6").Value = '''  -3.55%  '
$ws.Range("D7").Value = '''0.555'
$ws.Range("E7").Value = '''  -2.59%  '
$ws.Range("E8").Value = '''  +0.10%  '
$ws.Range("D9").Value = '''0.606'
$ws.Range("E9").Value = '''  -3.75%  '
$ws.Range("D10").Value = '''37.92'
$ws.Range("E10").Value = '''  -3.43%  '
$ws.Range("E11").Value = '''  +2.82%  '
$ws.Range("D12").Value = '''0.0854'
$ws.Range("E12").Value = '''  -3.38%  '
$ws.Range("E13").Value = '''  -3.91%  '
$ws.Range("D14").Value = '''3.467.88'
$ws.Range("E14").Value = '''  +2.46%  '
$ws.Range("D15").Value = '''7.57'
$ws.Range("E15").Value = '''  -4.05%  '
$ws.Range("D16").Value = '''2.980.20'
$ws.Range("E16").Value = '''  +1.95%  '
$ws.Range("E17").Value = '''  +1.96%  '
$ws.Range("D18").Value = '''51.936.42'
$ws.Range("E18").Value = '''  +0.25%  '
$ws.Range("E19").Value = '''  +2.10%  '
$ws.Range("D20").Value = '''7.43'
$ws.Range("E20").Value = '''  -1.42%  '
$ws.Range("D21").Value = '''13.50'
$ws.Range("E21").Value = '''  -3.77%  '
$ws.Range("D22").Value = '''0.0₃0968'
$ws.Range("E22").Value = '''  -1.49%  '
$ws.Range("D23").Value = '''68.92'
$ws.Range("E23").Value = '''  -2.70%  '
$ws.Range("D24").Value = '''262.96'
$ws.Range("E24").Value = '''  -2.81%  '
$ws.Range("E25").Value = '''  -3.89%  '
$ws.Range("D26").Value = '''0.178'
$ws.Range("E26").Value = '''  -2.04%  '
$ws.Range("D27").Value = '''26.87'
$ws.Range("E27").Value = '''  -0.96%  '
$ws.Range("E28").Value = '''  -0.01%  '
$ws.Range("D29").Value = '''7.38'
$ws.Range("E29").Value = '''  +0.89%  '
$ws.Range("E30").Value = '''  +1.87%  '
$ws.Range("E31").Value = '''  +4.95%  '
$ws.Range("D32").Value = '''10.13'
$ws.Range("E32").Value = '''  -4.16%  '
$ws.Range("D33").Value = '''35.98'
$ws.Range("E33").Value = '''  -7.49%  '
$ws.Range("E34").Value = '''  +12.99%  '
$ws.Range("D35").Value = '''51.24'
$ws.Range("E35").Value = '''  -1.70%  '
$ws.Range("E36").Value = '''  -2.09%  '
$ws.Range("E37").Value = '''  -0.05%  '
$ws.Range("D38").Value = '''3.29'
$ws.Range("E38").Value = '''  +1.73%  '
$ws.Range("D39").Value = '''2.83'
$ws.Range("E39").Value = '''  +1.77%  '
$ws.Range("D40").Value = '''17.47'
$ws.Range("E40").Value = '''  -5.17%  '
$ws.Range("D41").Value = '''1.93'
$ws.Range("E41").Value = '''  -3.59%  '
$ws.Range("D43").Value = '''23.06'
$ws.Range("E43").Value = '''  +0.11%  '
$ws.Range("D44").Value = '''123.22'
$ws.Range("E44").Value = '''  +2.99%  '
$ws.Range("D45").Value = '''2.17'
$ws.Range("E45").Value = '''  -0.27%  '
$ws.Range("D46").Value = '''2.124.08'
$ws.Range("E46").Value = '''  -0.72%  '
$ws.Range("E47").Value = '''  -4.47%  '
$ws.Range("E48").Value = '''  -7.40%  '
$ws.Range("B49").Value = '''TheGraph'
$ws.Range("C49").Value = '''https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D49").Value = '''0.242'
$ws.Range("E49").Value = '''  -3.67%  '
$ws.Range("B50").Value = '''BEAM'
$ws.Range("C50").Value = '''https://coinranking.com/coin/cYYMfXF4u+beam-beam'
$ws.Range("D50").Value = '''0.0334'
$ws.Range("E50").Value = '''  -1.15%  '
$ws.Range("B51").Value = '''SEI'
$ws.Range("C51").Value = '''https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D51").Value = '''0.900'
$ws.Range("E51").Value = '''  -0.87%  '
